$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "418.99" or "43.80" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.215.20"
$ws.Range("E2").Value = "  +4.66%  "
$ws.Range("D3").Value = "3.509.21"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "418.99"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "132.72"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  +4.27%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("D9").Value = "0.783"
$ws.Range("E9").Value = "  +7.39%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +15.58%  "
$ws.Range("D11").Value = "43.80"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "0.0000269"
$ws.Range("E12").Value = "  +22.34%  "
$ws.Range("D13").Value = "10.21"
$ws.Range("E13").Value = "  +10.44%  "
$ws.Range("D14").Value = "4.063.87"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "20.69"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "3.482.69"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "12.97"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").Value = "65.118.88"
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("D21").Value = "454.25"
$ws.Range("E21").Value = "  -5.15%  "
$ws.Range("D22").Value = "90.40"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "3.25"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "13.40"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").Value = "3.45"
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").Value = "10.03"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("D27").Value = "34.27"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("D28").Value = "12.76"
$ws.Range("E28").Value = "  +7.56%  "
$ws.Range("D29").Value = "2.74"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "7.51"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  +5.39%  "
$ws.Range("D32").Value = "0.163"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "40.05"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "57.58"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  +4.53%  "
$ws.Range("D37").Value = "0.0₃0750"
$ws.Range("E37").Value = "  +40.54%  "
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  +10.82%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "4.59"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").Value = "146.23"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "3.30"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").Value = "15.99"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").Value = "2.59"
$ws.Range("E50").Value = "  +12.33%  "
$ws.Range("D51").Value = "21.80"
$ws.Range("E51").Value = "  -1.70%  "

# Remove the temporary text formatting override so the cells
# retain their original (default) style, matching the source data.
$ws.Range("D2:D51").ClearFormats()
